# [DM] [delete field] scinario_12
#
# Remove the "skill_cd" sample row (row 16) and the "スキル" (skill)
# field column (column H) from the item master sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete the leftover "skill_cd" data row entirely; this shifts every
# row below it up by one.
$ws.Rows.Item(16).Delete() | Out-Null

# Delete the "スキル" field column entirely; this shifts every column
# to its right left by one (there is none here, H was the last used
# column).
$ws.Columns.Item(8).Delete() | Out-Null

# Restore the active selection to a sane cell now that the deleted
# row/column are gone.
$ws.Range("F7").Select() | Out-Null
